$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C2 to be a true numeric value (was stored as text "8")
$ws.Cells.Item(2, 3).Value = 8

# Ensure the date column keeps being stored as plain text, not real dates
$ws.Range("A3:A5").NumberFormat = "@"

# Row 3
$ws.Cells.Item(3, 1).Value = "2024-08-06"
$ws.Cells.Item(3, 2).Value = "sauju basnet"
$ws.Cells.Item(3, 3).Value = 7

# Row 4
$ws.Cells.Item(4, 1).Value = "2024-08-07"
$ws.Cells.Item(4, 2).Value = "sumit bam"
$ws.Cells.Item(4, 3).Value = 8

# Row 5 - Student-Semester stays as text "7"
$ws.Cells.Item(5, 1).Value = "2024-08-07"
$ws.Cells.Item(5, 2).Value = "sauju basnet"
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "7"
